$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Bmp7"
$ws.Range("C2").Value = "Acvr2b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.668521
$ws.Range("H2").Value = 5.005563
$ws.Range("I2").Value = 0.9677024783929865
$ws.Range("J2").Value = 0.9677024783929865
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.539665666666667
$ws.Range("N2").Value = 4.618997
$ws.Range("O2").Value = 0.3572088291809875
$ws.Range("P2").Value = 0.3572088291809875
$ws.Range("Q2").Value = 2.568964497812334
$ws.Range("R2").Value = 23.120680480311
$ws.Range("S2").Value = 0.3456718693022986
$ws.Range("T2").Value = 0.3456718693022986
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Bmp7"
$ws.Range("C3").Value = "Acvr2b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.668521
$ws.Range("H3").Value = 5.005563
$ws.Range("I3").Value = 0.9677024783929865
$ws.Range("J3").Value = 0.9677024783929865
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.452872333333333
$ws.Range("N3").Value = 4.358617
$ws.Range("O3").Value = 0.3370724153789985
$ws.Range("P3").Value = 0.3370724153789985
$ws.Range("Q3").Value = 2.424147998485667
$ws.Range("R3").Value = 21.817331986371
$ws.Range("S3").Value = 0.3261858117601671
$ws.Range("T3").Value = 0.3261858117601671
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("C4").Value = "Acvr2b"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.668521
$ws.Range("H4").Value = 5.005563
$ws.Range("I4").Value = 0.9677024783929865
$ws.Range("J4").Value = 0.9677024783929865
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.317729666666667
$ws.Range("N4").Value = 3.953189
$ws.Range("O4").Value = 0.3057187554400141
$ws.Range("P4").Value = 0.3057187554400141
$ws.Range("Q4").Value = 2.198659621156334
$ws.Range("R4").Value = 19.787936590407
$ws.Range("S4").Value = 0.2958447973305209
$ws.Range("T4").Value = 0.2958447973305209
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Bmp7"
$ws.Range("C5").Value = "Acvr2b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05568766666666666
$ws.Range("H5").Value = 0.167063
$ws.Range("I5").Value = 0.03229752160701353
$ws.Range("J5").Value = 0.03229752160701353
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.539665666666667
$ws.Range("N5").Value = 4.618997
$ws.Range("O5").Value = 0.3572088291809875
$ws.Range("P5").Value = 0.3572088291809875
$ws.Range("Q5").Value = 0.08574038842344445
$ws.Range("R5").Value = 0.771663495811
$ws.Range("S5").Value = 0.01153695987868895
$ws.Range("T5").Value = 0.01153695987868895
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Bmp7"
$ws.Range("C6").Value = "Acvr2b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05568766666666666
$ws.Range("H6").Value = 0.167063
$ws.Range("I6").Value = 0.03229752160701353
$ws.Range("J6").Value = 0.03229752160701353
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.452872333333333
$ws.Range("N6").Value = 4.358617
$ws.Range("O6").Value = 0.3370724153789985
$ws.Range("P6").Value = 0.3370724153789985
$ws.Range("Q6").Value = 0.08090707020788888
$ws.Range("R6").Value = 0.7281636318709999
$ws.Range("S6").Value = 0.01088660361883144
$ws.Range("T6").Value = 0.01088660361883144
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Bmp7"
$ws.Range("C7").Value = "Acvr2b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05568766666666666
$ws.Range("H7").Value = 0.167063
$ws.Range("I7").Value = 0.03229752160701353
$ws.Range("J7").Value = 0.03229752160701353
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.317729666666667
$ws.Range("N7").Value = 3.953189
$ws.Range("O7").Value = 0.3057187554400141
$ws.Range("P7").Value = 0.3057187554400141
$ws.Range("Q7").Value = 0.07338129043411111
$ws.Range("R7").Value = 0.660431613907
$ws.Range("S7").Value = 0.009873958109493141
$ws.Range("T7").Value = 0.009873958109493141